$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- copy the date-cell formatting (style used by column A) down into the new rows ---
$ws.Range("A48").Copy()
$ws.Range("A49:A51").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- data: note B is entered before A on each row so the dependent shared
#     formula in column E recalculates correctly ---
$ws.Range("B49").Value = 26000
$ws.Range("A49").Value = 43866
$ws.Range("D49").Value = "Ordered Amount"
$ws.Range("D49").Font.Italic = $false

$ws.Range("B50").Value = 45760
$ws.Range("A50").Value = 43867
$ws.Range("D50").Value = "Ordered Amount"

$ws.Range("B51").Value = 29120
$ws.Range("A51").Value = 43868
$ws.Range("D51").Value = "Ordered Amount"

# last populated row gets the italic "current row" styling
$ws.Range("D51").Font.Italic = $true

# --- view state: new active cell selection (pane stays frozen on row 1) ---
$null = $ws.Range("D56").Select()
